# Apply the "major edits to crypto full course" change:
#  1. Update cached date-field text from 2/6/2023 -> 2/9/2024 everywhere it
#     appears (slide master, all slide layouts, and the notes master).
#  2. Expand a bullet on the "AES Extras (2)" slide to mention "elements".

$p = $ppt.ActivePresentation

$oldDate = "2/6/2023"
$newDate = "2/9/2024"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master date placeholder.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# Slide 15 ("AES Extras (2)") bullet text update.
$slide = $p.Slides.Item(15)
$contentShape = $slide.Shapes.Item(2)
$textRange = $contentShape.TextFrame.TextRange
$oldBullet = "Redefines addition, multiplication, and inverse"
$newBullet = "Redefines elements, addition, multiplication, and inverse"
$paraCount = $textRange.Paragraphs().Count
for ($pi = 1; $pi -le $paraCount; $pi++) {
    $para = $textRange.Paragraphs($pi, 1)
    # Paragraphs(...).Text includes the trailing paragraph-mark (`r),
    # so trim it off before comparing against the plain bullet text.
    if ($para.Text.TrimEnd("`r") -eq $oldBullet) {
        $para.Text = $newBullet
    }
}
